$d = $word.ActiveDocument

# Locate the "la) Crea un baseviewmodel que implemente on property." text
# (spans several runs) and mark it in red, matching the other lettered
# exercise items' highlight color.
$rng = $d.Content
$found = $rng.Find.Execute("la) Crea un baseviewmodel que implemente on property.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Font.Color = 255  # wdColorRed (0x0000FF -> BGR = FF0000 in OOXML)

    # Move the "_GoBack" bookmark (tracking the last edit location) to sit
    # right after the text we just edited. Re-adding a bookmark with the
    # same name as an existing one moves it, removing the previous one.
    $endRng = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $endRng)
}
